# Updates cryptos list cell values (price + 1h volume %) per the Feb 27 2024 GitHub Actions refresh.
# Rows 41/42 (Monero / Celestia) also swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.637.29"
$ws.Range("E2").Value = "  +11.05%  "
# Row 3
$ws.Range("D3").Value = "3.249.19"
$ws.Range("E3").Value = "  +6.77%  "
# Row 4
$ws.Range("E4").Value = "  -0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.54%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.51%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.53%  "
# Row 8
$ws.Range("E8").Value = "  -0.08%  "
# Row 9
$ws.Range("E9").Value = "  +8.29%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.46"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.37%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +12.41%  "
# Row 12
$ws.Range("E12").Value = "  +2.58%  "
# Row 13
$ws.Range("D13").Value = "3.759.22"
$ws.Range("E13").Value = "  +6.22%  "
# Row 14
$ws.Range("E14").Value = "  +6.53%  "
# Row 15
$ws.Range("E15").Value = "  +5.48%  "
# Row 16
$ws.Range("D16").Value = "3.254.10"
$ws.Range("E16").Value = "  +4.63%  "
# Row 17
$ws.Range("E17").Value = "  +7.31%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.00%  "
# Row 19
$ws.Range("D19").Value = "56.481.61"
$ws.Range("E19").Value = "  +10.66%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.25%  "
# Row 21
$ws.Range("E21").Value = "  +9.78%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.98%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "300.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +14.21%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.79"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.22%  "
# Row 25
$ws.Range("E25").Value = "  +4.36%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.72%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.26%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.35"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.19%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.79%  "
# Row 30
$ws.Range("E30").Value = "  +5.45%  "
# Row 31
$ws.Range("E31").Value = "  +0.04%  "
# Row 32
$ws.Range("E32").Value = "  +6.87%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.74%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.60%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0490"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.19%  "
# Row 36
$ws.Range("E36").Value = "  +3.62%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.92%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.17%  "
# Row 39
$ws.Range("E39").Value = "  -0.18%  "
# Row 40
$ws.Range("E40").Value = "  +26.28%  "
# Row 41
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.88%  "
# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.03"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.46%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.57%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.11%  "
# Row 45
$ws.Range("E45").Value = "  +5.15%  "
# Row 46
$ws.Range("E46").Value = "  -0.68%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.10%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +56.99%  "
# Row 49
$ws.Range("D49").Value = "2.143.45"
$ws.Range("E49").Value = "  +4.60%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.44%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.43"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.46%  "
